# Auto-generated Excel COM-interop script applying the scheduled market-data refresh
# described by the commit "chore: update Sheets via scheduled runner".
# For every affected sheet/row, rewrite the updated H:N crafting-profit columns.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 588.8570999999999
$ws.Range("I12").Value = 484.4
$ws.Range("J12").Value = 850
$ws.Range("K12").Value = 484.4
$ws.Range("L12").Value = 850
$ws.Range("M12").Value = -314.4
$ws.Range("N12").Value = -1190

$ws.Range("H74").Value = 3745.8
$ws.Range("I74").Value = 3780.6365
$ws.Range("K74").Value = 3780.6365
$ws.Range("M74").Value = -2844.6365

$ws.Range("H77").Value = 3745.8
$ws.Range("I77").Value = 3780.6365
$ws.Range("K77").Value = 18903.1825
$ws.Range("M77").Value = -14223.1825

$ws.Range("H131").Value = 3421.7917
$ws.Range("I131").Value = 1915.3334
$ws.Range("J131").Value = 5932.5557
$ws.Range("K131").Value = 5746.0002
$ws.Range("L131").Value = 17797.6671
$ws.Range("M131").Value = -706.0002000000004
$ws.Range("N131").Value = -27877.6671

$ws.Range("H137").Value = 389517.44
$ws.Range("I137").Value = 626440.25
$ws.Range("K137").Value = 1879320.75
$ws.Range("M137").Value = -1876770.75

$ws.Range("H138").Value = 2327.46
$ws.Range("I138").Value = 994.2727
$ws.Range("J138").Value = 2703.487
$ws.Range("K138").Value = 2982.8181
$ws.Range("L138").Value = 8110.461
$ws.Range("M138").Value = 2157.1819
$ws.Range("N138").Value = -18390.461

$ws.Range("H140").Value = 78956
$ws.Range("J140").Value = 78956
$ws.Range("L140").Value = 78956
$ws.Range("N140").Value = -89316

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 34094880
$ws.Range("I61").Value = 26318690
$ws.Range("K61").Value = 26318690
$ws.Range("M61").Value = -26318478

$ws.Range("H74").Value = 6256565
$ws.Range("I74").Value = 8930511
$ws.Range("K74").Value = 8930511
$ws.Range("M74").Value = -8929637

$ws.Range("H77").Value = 6256565
$ws.Range("I77").Value = 8930511
$ws.Range("K77").Value = 44652555
$ws.Range("M77").Value = -44648187

$ws.Range("H81").Value = 80870
$ws.Range("J81").Value = 97830
$ws.Range("L81").Value = 97830
$ws.Range("N81").Value = -99826

$ws.Range("H84").Value = 80870
$ws.Range("J84").Value = 97830
$ws.Range("L84").Value = 293490
$ws.Range("N84").Value = -303474

$ws.Range("H124").Value = 57495
$ws.Range("J124").Value = 57495
$ws.Range("L124").Value = 57495
$ws.Range("N124").Value = -67315

$ws.Range("H125").Value = 52996.332
$ws.Range("J125").Value = 52996.332
$ws.Range("L125").Value = 52996.332
$ws.Range("N125").Value = -62836.332

$ws.Range("H132").Value = 9263923
$ws.Range("I132").Value = 11906877
$ws.Range("K132").Value = 35720631
$ws.Range("M132").Value = -35718101

$ws.Range("H136").Value = 34094880
$ws.Range("I136").Value = 26318690
$ws.Range("K136").Value = 78956070
$ws.Range("M136").Value = -78953520

$ws.Range("H139").Value = 75949.25
$ws.Range("J139").Value = 75949.25
$ws.Range("L139").Value = 75949.25
$ws.Range("N139").Value = -86229.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 39473.8
$ws.Range("I26").Value = 15790
$ws.Range("J26").Value = 74999.5
$ws.Range("K26").Value = 15790
$ws.Range("L26").Value = 74999.5
$ws.Range("M26").Value = -15498
$ws.Range("N26").Value = -75583.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 101322.664
$ws.Range("J20").Value = 101322.664
$ws.Range("L20").Value = 101322.664
$ws.Range("N20").Value = -101794.664

$ws.Range("H30").Value = 101322.664
$ws.Range("J30").Value = 101322.664
$ws.Range("L30").Value = 101322.664
$ws.Range("N30").Value = -101504.664

$ws.Range("H31").Value = 1311723.1
$ws.Range("I31").Value = 26582.5
$ws.Range("J31").Value = 1954293.4
$ws.Range("K31").Value = 26582.5
$ws.Range("L31").Value = 1954293.4
$ws.Range("M31").Value = -26287.5
$ws.Range("N31").Value = -1954883.4

$ws.Range("H34").Value = 1311723.1
$ws.Range("I34").Value = 26582.5
$ws.Range("J34").Value = 1954293.4
$ws.Range("K34").Value = 26582.5
$ws.Range("L34").Value = 1954293.4
$ws.Range("M34").Value = -26380.5
$ws.Range("N34").Value = -1954697.4

$ws.Range("H58").Value = 3703.3125
$ws.Range("J58").Value = 3544.111
$ws.Range("L58").Value = 3544.111
$ws.Range("N58").Value = -3950.111

$ws.Range("H105").Value = 1647.2667
$ws.Range("J105").Value = 1058
$ws.Range("L105").Value = 1058
$ws.Range("N105").Value = -4552

$ws.Range("H108").Value = 75962
$ws.Range("J108").Value = 77194.39999999999
$ws.Range("L108").Value = 77194.39999999999
$ws.Range("N108").Value = -84874.39999999999

$ws.Range("H119").Value = 65495
$ws.Range("J119").Value = 65495
$ws.Range("L119").Value = 65495
$ws.Range("N119").Value = -75171

$ws.Range("H122").Value = 2549.4736
$ws.Range("J122").Value = 2041.8889
$ws.Range("L122").Value = 6125.6667
$ws.Range("N122").Value = -11025.6667

$ws.Range("H128").Value = 101322.664
$ws.Range("J128").Value = 101322.664
$ws.Range("L128").Value = 101322.664
$ws.Range("N128").Value = -111282.664

$ws.Range("H136").Value = 3703.3125
$ws.Range("J136").Value = 3544.111
$ws.Range("L136").Value = 10632.333
$ws.Range("N136").Value = -15732.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 854.8
$ws.Range("I121").Value = 64.166664
$ws.Range("J121").Value = 1193.6428
$ws.Range("K121").Value = 192.499992
$ws.Range("L121").Value = 3580.9284
$ws.Range("M121").Value = 1117.500008
$ws.Range("N121").Value = -6200.928400000001

$ws.Range("H131").Value = 4558.5293
$ws.Range("I131").Value = 9589.166999999999
$ws.Range("J131").Value = 1814.5454
$ws.Range("K131").Value = 28767.501
$ws.Range("L131").Value = 5443.6362
$ws.Range("M131").Value = -23727.501
$ws.Range("N131").Value = -15523.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 446.81818
$ws.Range("I107").Value = 202.14285
$ws.Range("K107").Value = 202.14285
$ws.Range("M107").Value = 1717.85715

$ws.Range("H126").Value = 4751.5
$ws.Range("I126").Value = 3506
$ws.Range("J126").Value = 5166.6665
$ws.Range("K126").Value = 10518
$ws.Range("L126").Value = 15499.9995
$ws.Range("M126").Value = -8048
$ws.Range("N126").Value = -20439.9995

$ws.Range("H136").Value = 6679.4346
$ws.Range("J136").Value = 6679.4346
$ws.Range("L136").Value = 20038.3038
$ws.Range("N136").Value = -25138.3038

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 486455.3
$ws.Range("I7").Value = 7217.125
$ws.Range("K7").Value = 7217.125
$ws.Range("M7").Value = -7105.125

$ws.Range("H40").Value = 4107.5454
$ws.Range("I40").Value = 3557.9333
$ws.Range("K40").Value = 3557.9333
$ws.Range("M40").Value = -3421.9333

$ws.Range("H100").Value = 2469.8
$ws.Range("J100").Value = 3300
$ws.Range("L100").Value = 3300
$ws.Range("N100").Value = -4382

$ws.Range("H126").Value = 486455.3
$ws.Range("I126").Value = 7217.125
$ws.Range("K126").Value = 21651.375
$ws.Range("M126").Value = -19181.375

$ws.Range("H133").Value = 72000
$ws.Range("J133").Value = 72000
$ws.Range("L133").Value = 72000
$ws.Range("N133").Value = -77060

$ws.Range("H136").Value = 34241.676
$ws.Range("I136").Value = 6040.0835
$ws.Range("K136").Value = 18120.2505
$ws.Range("M136").Value = -15570.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3251.7827
$ws.Range("I122").Value = 3339.8635
$ws.Range("J122").Value = 1314
$ws.Range("K122").Value = 10019.5905
$ws.Range("L122").Value = 3942
$ws.Range("M122").Value = -7569.5905
$ws.Range("N122").Value = -8842

$ws.Range("H130").Value = 132000
$ws.Range("J130").Value = 132000
$ws.Range("L130").Value = 132000
$ws.Range("N130").Value = -142040

$ws.Range("H136").Value = 2797.2856
$ws.Range("I136").Value = 2484.2563
$ws.Range("K136").Value = 7452.7689
$ws.Range("M136").Value = -4902.7689
